$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# BagItem_背包物品表 - rows 9-28 hold 4 groups of 5 test rows, one group per
# "DragonName0000N" placeholder name, with sequential ids (4..23), a shared
# icon id (108377), achievable = 1 and a shared description
# ("TestBagItemDesc0001"). Fix up the cell formatting first (copying the
# already-correct style patterns established by the existing rows) and then
# (re)write every value so the final content matches regardless of paste
# order.
# ---------------------------------------------------------------------------

# --- Column A: wrap/no-wrap style cadence is 1,5,5,1,1 per 5-row group ------
$ws.Range("A10").Copy()
$ws.Range($ws.Cells.Item(10,1), $ws.Cells.Item(11,1)).PasteSpecial(-4122)
$ws.Range("A9").Copy()
$ws.Range($ws.Cells.Item(12,1), $ws.Cells.Item(14,1)).PasteSpecial(-4122)
$ws.Range("A10").Copy()
$ws.Range($ws.Cells.Item(15,1), $ws.Cells.Item(16,1)).PasteSpecial(-4122)
$ws.Range("A9").Copy()
$ws.Range($ws.Cells.Item(17,1), $ws.Cells.Item(19,1)).PasteSpecial(-4122)
$ws.Range("A10").Copy()
$ws.Range($ws.Cells.Item(20,1), $ws.Cells.Item(21,1)).PasteSpecial(-4122)
$ws.Range("A9").Copy()
$ws.Range($ws.Cells.Item(22,1), $ws.Cells.Item(24,1)).PasteSpecial(-4122)
$ws.Range("A10").Copy()
$ws.Range($ws.Cells.Item(25,1), $ws.Cells.Item(26,1)).PasteSpecial(-4122)
$ws.Range("A9").Copy()
$ws.Range($ws.Cells.Item(27,1), $ws.Cells.Item(28,1)).PasteSpecial(-4122)

# --- Columns B-E: constant style for the whole 9-28 block -------------------
$ws.Range("B9").Copy()
$ws.Range($ws.Cells.Item(10,2), $ws.Cells.Item(28,2)).PasteSpecial(-4122)
$ws.Range("C9").Copy()
$ws.Range($ws.Cells.Item(10,3), $ws.Cells.Item(28,3)).PasteSpecial(-4122)
$ws.Range("D9").Copy()
$ws.Range($ws.Cells.Item(10,4), $ws.Cells.Item(28,4)).PasteSpecial(-4122)
$ws.Range("E9").Copy()
$ws.Range($ws.Cells.Item(10,5), $ws.Cells.Item(28,5)).PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Now (re)write the actual values ----------------------------------------
$names = @("DragonName00001", "DragonName00002", "DragonName00003", "DragonName00004")

$id = 4
$row = 9
foreach ($name in $names) {
    for ($i = 0; $i -lt 5; $i++) {
        $ws.Cells.Item($row, 1).Value = $id
        $ws.Cells.Item($row, 2).Value = $name
        $ws.Cells.Item($row, 3).Value = "TestBagItemDesc0001"
        $ws.Cells.Item($row, 4).Value = 108377
        $ws.Cells.Item($row, 5).Value = 1
        $id = $id + 1
        $row = $row + 1
    }
}

# ---------------------------------------------------------------------------
# Selection moved from C13 to C9
# ---------------------------------------------------------------------------
$ws.Range("C9").Select()
